$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update Mid X / Mid Y values for the five *A designators (C1A, D1A, R1A, R2A, R3A)
$ws.Range("B13").Value = "0.470000"
$ws.Range("C13").Value = "42.120000"

$ws.Range("B31").Value = "-4.760000"
$ws.Range("C31").Value = "41.890000"

$ws.Range("B49").Value = "-1.980000"
$ws.Range("C49").Value = "42.110000"

$ws.Range("B57").Value = "2.790000"
$ws.Range("C57").Value = "41.980000"

$ws.Range("B59").Value = "5.160000"
$ws.Range("C59").Value = "42.190000"

# Update the selection to A2:E68 with active cell A2
$ws.Range("A2:E68").Select()
